$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 347
$ws.Range("F2").Value2 = 45597
$ws.Range("G2").Value2 = 30468
$ws.Range("H2").Value2 = 45658
$ws.Range("E3").Value2 = 30407
$ws.Range("F3").Value2 = 45597
$ws.Range("E4").Value2 = 30407
$ws.Range("F4").Value2 = 45597
$ws.Range("G4").Value2 = 30468
$ws.Range("H4").Value2 = 45658
$ws.Range("E5").Value2 = 30407
$ws.Range("F5").Value2 = 45597
$ws.Range("G5").Value2 = 30468
$ws.Range("H5").Value2 = 45658
$ws.Range("C6").Value2 = 455
$ws.Range("F6").Value2 = 45566
$ws.Range("G6").Value2 = 30468
$ws.Range("H6").Value2 = 45658
$ws.Range("E7").Value2 = 30376
$ws.Range("F7").Value2 = 45566
$ws.Range("G7").Value2 = 30468
$ws.Range("H7").Value2 = 45658
$ws.Range("D8").Value2 = 423
$ws.Range("E8").Value2 = 30376
$ws.Range("F8").Value2 = 45566
$ws.Range("H8").Value2 = 45658
$ws.Range("E9").Value2 = 30407
$ws.Range("F9").Value2 = 45597
$ws.Range("G9").Value2 = 30468
$ws.Range("H9").Value2 = 45658
$ws.Range("D10").Value2 = 496
$ws.Range("E10").Value2 = 30407
$ws.Range("F10").Value2 = 45597
$ws.Range("H10").Value2 = 45658
$ws.Range("E11").Value2 = 30376
$ws.Range("F11").Value2 = 45566
$ws.Range("G11").Value2 = 30468
$ws.Range("H11").Value2 = 45658
$ws.Range("C12").Value2 = 384
$ws.Range("D12").Value2 = 365
$ws.Range("F12").Value2 = 45597
$ws.Range("H12").Value2 = 45658
$ws.Range("C13").Value2 = 480
$ws.Range("F13").Value2 = 45597
$ws.Range("G13").Value2 = 30468
$ws.Range("H13").Value2 = 45658
$ws.Range("C14").Value2 = 435
$ws.Range("D14").Value2 = 409
$ws.Range("F14").Value2 = 45566
$ws.Range("H14").Value2 = 45658
$ws.Range("G15").Value2 = 30438
$ws.Range("H15").Value2 = 45658
$ws.Range("C16").Value2 = 467
$ws.Range("D16").Value2 = 423
$ws.Range("F16").Value2 = 45566
$ws.Range("H16").Value2 = 45658
$ws.Range("C17").Value2 = 383
$ws.Range("D17").Value2 = 407
$ws.Range("F17").Value2 = 45597
$ws.Range("H17").Value2 = 45658
$ws.Range("D18").Value2 = 273
$ws.Range("E18").Value2 = 30407
$ws.Range("F18").Value2 = 45597
$ws.Range("H18").Value2 = 45658
$ws.Range("D19").Value2 = 411
$ws.Range("E19").Value2 = 30407
$ws.Range("F19").Value2 = 45597
$ws.Range("H19").Value2 = 45658
$ws.Range("C20").Value2 = 493
$ws.Range("F20").Value2 = 45566
$ws.Range("G20").Value2 = 30468
$ws.Range("H20").Value2 = 45658
$ws.Range("C21").Value2 = 323
$ws.Range("F21").Value2 = 45597
$ws.Range("G21").Value2 = 30468
$ws.Range("H21").Value2 = 45658
$ws.Range("C22").Value2 = 336
$ws.Range("D22").Value2 = 380
$ws.Range("F22").Value2 = 45597
$ws.Range("H22").Value2 = 45658
$ws.Range("C23").Value2 = 272
$ws.Range("D23").Value2 = 415
$ws.Range("F23").Value2 = 45536
$ws.Range("H23").Value2 = 45658
$ws.Range("D24").Value2 = 423
$ws.Range("H24").Value2 = 45658
$ws.Range("D25").Value2 = 328
$ws.Range("E25").Value2 = 30376
$ws.Range("F25").Value2 = 45566
$ws.Range("H25").Value2 = 45658
$ws.Range("C26").Value2 = 345
$ws.Range("D26").Value2 = 326
$ws.Range("F26").Value2 = 45566
$ws.Range("H26").Value2 = 45658
$ws.Range("G27").Value2 = 30468
$ws.Range("H27").Value2 = 45658
$ws.Range("D28").Value2 = 393
$ws.Range("E28").Value2 = 30376
$ws.Range("F28").Value2 = 45566
$ws.Range("H28").Value2 = 45658
$ws.Range("C29").Value2 = 274
$ws.Range("D29").Value2 = 249
$ws.Range("F29").Value2 = 45566
$ws.Range("H29").Value2 = 45658
$ws.Range("D30").Value2 = 230
$ws.Range("E30").Value2 = 30407
$ws.Range("F30").Value2 = 45597
$ws.Range("H30").Value2 = 45658
$ws.Range("C31").Value2 = 406
$ws.Range("F31").Value2 = 45566
$ws.Range("G31").Value2 = 30468
$ws.Range("H31").Value2 = 45658
$ws.Range("E32").Value2 = 30407
$ws.Range("F32").Value2 = 45597
$ws.Range("G32").Value2 = 30468
$ws.Range("H32").Value2 = 45658
$ws.Range("D33").Value2 = 411
$ws.Range("H33").Value2 = 45658
$ws.Range("C34").Value2 = 214
$ws.Range("D34").Value2 = 330
$ws.Range("F34").Value2 = 45566
$ws.Range("H34").Value2 = 45658
$ws.Range("C35").Value2 = 417
$ws.Range("D35").Value2 = 330
$ws.Range("F35").Value2 = 45566
$ws.Range("H35").Value2 = 45658
$ws.Range("D36").Value2 = 423
$ws.Range("E36").Value2 = 30407
$ws.Range("F36").Value2 = 45597
$ws.Range("H36").Value2 = 45658
$ws.Range("C37").Value2 = 479
$ws.Range("D37").Value2 = 330
$ws.Range("F37").Value2 = 45597
$ws.Range("H37").Value2 = 45658
$ws.Range("C38").Value2 = 372
$ws.Range("D38").Value2 = 380
$ws.Range("F38").Value2 = 45597
$ws.Range("H38").Value2 = 45658
$ws.Range("C39").Value2 = 240
$ws.Range("D39").Value2 = 236
$ws.Range("F39").Value2 = 45597
$ws.Range("H39").Value2 = 45658
$ws.Range("C40").Value2 = 299
$ws.Range("D40").Value2 = 327
$ws.Range("F40").Value2 = 45597
$ws.Range("H40").Value2 = 45658
$ws.Range("C41").Value2 = 406
$ws.Range("D41").Value2 = 249
$ws.Range("F41").Value2 = 45566
$ws.Range("H41").Value2 = 45658
$ws.Range("C42").Value2 = 249
$ws.Range("D42").Value2 = 231
$ws.Range("F42").Value2 = 45566
$ws.Range("H42").Value2 = 45658
$ws.Range("C43").Value2 = 485
$ws.Range("D43").Value2 = 330
$ws.Range("F43").Value2 = 45566
$ws.Range("H43").Value2 = 45658
$ws.Range("C44").Value2 = 418
$ws.Range("D44").Value2 = 317
$ws.Range("F44").Value2 = 45566
$ws.Range("H44").Value2 = 45658
$ws.Range("D45").Value2 = 330
$ws.Range("H45").Value2 = 45658
$ws.Range("C46").Value2 = 346
$ws.Range("D46").Value2 = 311
$ws.Range("F46").Value2 = 45566
$ws.Range("H46").Value2 = 45658
$ws.Range("C47").Value2 = 348
$ws.Range("D47").Value2 = 273
$ws.Range("F47").Value2 = 45597
$ws.Range("H47").Value2 = 45658
$ws.Range("C48").Value2 = 371
$ws.Range("D48").Value2 = 328
$ws.Range("F48").Value2 = 45597
$ws.Range("H48").Value2 = 45658
$ws.Range("C49").Value2 = 310
$ws.Range("D49").Value2 = 326
$ws.Range("F49").Value2 = 45566
$ws.Range("H49").Value2 = 45658
$ws.Range("C50").Value2 = 371
$ws.Range("D50").Value2 = 251
$ws.Range("F50").Value2 = 45597
$ws.Range("H50").Value2 = 45658
$ws.Range("D51").Value2 = 330
$ws.Range("H51").Value2 = 45658
$ws.Range("C52").Value2 = 354
$ws.Range("D52").Value2 = 328
$ws.Range("F52").Value2 = 45566
$ws.Range("H52").Value2 = 45658
